$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_ETR"

$values = @(
    0,
    43.3690532265175,
    40.6347900735476,
    5.31446236928425,
    0,
    29.0978355337041,
    28.9485701599378,
    11.5345010103731,
    0,
    38.3432440967571,
    37.9321949419935,
    0,
    0,
    34.3865741783422,
    34.1169566734631,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
